$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Surveys")
$ws.Range("B12").Copy()
$ws.Range("A22:L25").PasteSpecial(-4122)
